$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph right after the
# table. We keep that paragraph (it becomes the first of two blank spacer
# paragraphs) and append four more paragraphs after it:
#   - a blank spacer paragraph
#   - the long "tell me about yourself" answer paragraph
#   - a blank spacer paragraph
#   - the closing sentence paragraph
$lastP = $d.Paragraphs.Last
$r = $lastP.Range

$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$bodyParaIndex = $n - 2
$closingParaIndex = $n

$bodyText = "I’m software engineer working for Oracle for 4 years. I’m working on the product called Snapshot Manager (SMU) which is the Java web application and service that provides the data protection solutions for oracle databases. I have two major roles or responsibilities: one is core product feature development mostly backend and API. And the other is to fix critical customer issue and design features and integration solutions based on customer request. "
$closingText = "I joined the team when the product is in the beta phase and now it becomes more mature."

$d.Paragraphs.Item($bodyParaIndex).Range.Text = $bodyText
$d.Paragraphs.Item($closingParaIndex).Range.Text = $closingText
